# feat: add 2022-Q4 data
#
# A new quarterly sheet "2022-Q4" is inserted right after "总计", pushing
# 2022-Q3 / 2022-Q2 / 2022-Q1 / 2021-Q4 one position to the right (their
# own data is untouched). The "总计" (totals) sheet gains a new top data
# row for 2022-Q4 and keeps its existing rows (now one row longer).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update the "总计" (totals) sheet: shift the four existing data rows
#    down by one (2022-Q3->row3, 2022-Q2->row4, 2022-Q1->row5,
#    2021-Q4->row6) and put the new 2022-Q4 figures in row 2.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Grab the style-carrying column-A cell so the brand new row 6 matches
# the look of the existing index cells (s="2" bold/centered/bordered).
$total.Range("A5").Copy($total.Range("A6"))

$total.Range("A6").Value = 4
$total.Range("B6").Value = "2021-Q4"
$total.Range("C6").Value = 2
$total.Range("D6").Value = 0.67

$total.Range("B5").Value = "2022-Q1"
$total.Range("D5").Value = 2.88

$total.Range("B4").Value = "2022-Q2"
$total.Range("D4").Value = 3.31

$total.Range("B3").Value = "2022-Q3"
$total.Range("D3").Value = 4.28

$total.Range("B2").Value = "2022-Q4"
$total.Range("D2").Value = 3.81

# ---------------------------------------------------------------------
# 2) Insert a brand new "2022-Q4" worksheet right after "总计", cloning
#    the layout/formatting of the existing "2022-Q3" sheet, then stamp
#    in the Q4-specific fund figures.
# ---------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("总计")
$wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet) | Out-Null

# Re-resolve by position (index 2 == right after "总计") instead of
# keeping the just-returned reference, which doesn't reliably act as a
# copy destination straight after creation.
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.UsedRange.Copy($q4.Range("A1"))

# The source sheet stores fund code / size / position figures as plain
# text (so "012348" keeps its leading zero, "39.65" isn't coerced to a
# number). Mark the destination range as Text first so COM's automatic
# type-sniffing on assignment leaves the strings alone.
$q4.Range("B2:G3").NumberFormat = "@"

$q4.Range("B2").Value = "012348"
$q4.Range("C2").Value = "天弘恒生科技指数（QDII）A"
$q4.Range("D2").Value = "39.65"
$q4.Range("E2").Value = "93.67"
$q4.Range("F2").Value = "4.94"
$q4.Range("G2").Value = "1.9587"
$q4.Range("H2").Value = 7

$q4.Range("B3").Value = "012349"
$q4.Range("C3").Value = "天弘恒生科技指数（QDII）C"
$q4.Range("D3").Value = "37.52"
$q4.Range("E3").Value = "93.67"
$q4.Range("F3").Value = "4.94"
$q4.Range("G3").Value = "1.8535"
$q4.Range("H3").Value = 7
